# Insert a new weekly record as row 313, pushing the existing "Apio"
# records (old rows 313-430) down by one row (to 314-431).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(313).Insert()

$ws.Range("A313").Value = 4
$ws.Range("B313").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C313").Value = "Los Lagos"
$ws.Range("D313").Value = 45009
$ws.Range("E313").Value = 10
$ws.Range("F313").Value = 100112017
$ws.Range("G313").Value = "Apio"
$ws.Range("H313").Value = "Americana (o)"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 50
$ws.Range("K313").Value = 11000
$ws.Range("L313").Value = 12000
$ws.Range("M313").Value = 11500
$ws.Range("N313").Value = "$/docena de matas"
$ws.Range("O313").Value = "Región de Coquimbo"
$ws.Range("P313").Value = 1917
$ws.Range("Q313").Value = 6
$ws.Range("R313").Value = "Hortaliza"
